$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.775.53"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.895.46"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'239.30"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4896"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").Value = "'0.2975"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").Value = "'0.06764"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "1.875.74"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "'16.96"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "'0.07306"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "'89.98"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "'5.128"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "'0.6694"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "30.725.88"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "'0.000007958"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'13.49"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "2.108.91"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "'0.9990"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "'211.75"
$ws.Range("E22").Value = "  +8.73%  "
$ws.Range("D23").Value = "'4.990"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("D24").Value = "'6.229"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("D25").Value = "'9.631"
$ws.Range("E25").Value = "  +2.99%  "
$ws.Range("D26").Value = "'158.60"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "'18.91"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").Value = "'1.894"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("D29").Value = "'1.418"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "'4.348"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'0.09155"
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").Value = "'4.054"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "'0.05186"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "'0.7513"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "'1.117"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").Value = "'2.688"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "'0.01835"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("D38").Value = "'2.694"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "'0.9278"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "'2.107"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").Value = "'0.4503"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("D42").Value = "'106.71"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'5.810"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D45").Value = "'7.822"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "'0.1376"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("D47").Value = "'65.72"
$ws.Range("E47").Value = "  +13.21%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4079"
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.92"
$ws.Range("E49").Value = "  +3.84%  "
$ws.Range("D50").Value = "'0.05896"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "'8.883"
$ws.Range("E51").Value = "  +1.25%  "
